# 4 Offizlieder + Fortschritt update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rows getting the "Offizium" marker in column B (with "ok" in column C)
$offiziumRows = @(21, 33, 37, 77)
foreach ($r in $offiziumRows) {
    $ws.Cells.Item($r, 2).Value = "Offizium"
    $ws.Cells.Item($r, 3).Value = "ok"
}

# Rows getting the "Lhymnen" marker in column B (with "ok" in column C)
$lhymnenRows = @(23, 25, 38, 39, 54, 55, 61, 89)
foreach ($r in $lhymnenRows) {
    $ws.Cells.Item($r, 2).Value = "Lhymnen"
    $ws.Cells.Item($r, 3).Value = "ok"
}

# Update the "Fortschritt" (progress) of the view: move the active
# selection to B78 (this also clears the previous topLeftCell scroll
# position that was left over at A53)
$ws.Range("B78").Select()
